$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of the existing header cell (H1) onto the new header cells
# so I1/J1 match the other header cells (bold, bordered, centered).
$headerStyle = $ws.Range("H1").Style
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1").Style = $headerStyle
$ws.Range("J1").Style = $headerStyle

# Fill in the new data columns I and J for rows 2-8
$values = @(
    @(7, 8),
    @(6, 8),
    @(8, 9),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(7, 7)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
